$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    'compression for stomach women',
    'compression for the knee',
    'compression for women',
    'compression full leg',
    'compression garment',
    'compression garment for women',
    'compression garments for legs',
    'compression garments for women',
    'compression garments legs',
    'compression garments women',
    'compression gear',
    'compression gear women',
    'compression gym leggings',
    'compression gym pants women',
    'compression hamstring',
    'compression high waist leggings women',
    'compression high waist yoga pants',
    'compression hip',
    'compression hip support',
    'compression hoes for women',
    'compression hose knee',
    'compression hose knee highs',
    'compression hose ladies',
    'compression jacket',
    'compression jacket women',
    'compression jackets',
    'compression joints',
    'compression knee',
    'compression knee high for women',
    'compression knee high women',
    'compression knee highs',
    'compression knee highs for women',
    'compression knee highs for women plus size',
    'compression knee highs plus size',
    'compression knee highs women',
    'compression knee leggings',
    'compression knee length',
    'compression knee pants',
    'compression knee plus size',
    'compression knee running',
    'compression knee support',
    'compression knee women',
    'compression knees',
    'compression leg',
    'compression leg band',
    'compression leg device',
    'compression leg support',
    'compression leg warmers',
    'compression legging',
    'compression leggings',
    'compression leggings 30',
    'compression leggings basketball',
    'compression leggings black',
    'compression leggings capri',
    'compression leggings circulation',
    'compression leggings for boys',
    'compression leggings for circulation',
    'compression leggings for plus size',
    'compression leggings for women',
    'compression leggings for women circulation',
    'compression leggings for women workout',
    'compression leggings high waist',
    'compression leggings low rise',
    'compression leggings medical',
    'compression leggings men basketball',
    'compression leggings men white',
    'compression leggings navy',
    'compression leggings plus',
    'compression leggings plus size',
    'compression leggings plus size women',
    'compression leggings pockets',
    'compression leggings running',
    'compression leggings running women',
    'compression leggings tall women',
    'compression leggings white',
    'compression leggings with pockets women',
    'compression leggings with side pocket',
    'compression leggings women',
    'compression leggings women 2xu',
    'compression leggings women active',
    'compression leggings women black',
    'compression leggings women circulation',
    'compression leggings women cwx',
    'compression leggings women for varicose veins',
    'compression leggings women high waist',
    'compression leggings women medical',
    'compression leggings women plus size',
    'compression leggings women running',
    'compression leggings women tall',
    'compression leggings xs',
    'compression leggins',
    'compression leggins women',
    'compression long shorts women',
    'compression long tank',
    'compression lower calf',
    'compression lycra',
    'compression medical leggings',
    'compression men 3 4',
    'compression men capri',
    'compression men leggings'
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
